# CoinGalaxy.docx edit
#
# 1. In the "Calander(not fix) - FIND THE API PF THIS" paragraph, mark a
#    grammar-check span: <w:proofErr w:type="gramStart"/> right after the
#    existing spellStart marker, split the "(not fix)" run into "(" and
#    "not fix)" runs, and close the grammar span with
#    <w:proofErr w:type="gramEnd"/> between them.
# 2. Bold the "Portfolio" paragraph (paragraph + run properties).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Find-ParagraphByPrefix($doc, [string]$prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like ($prefix + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Step 1: "Calander(not fix) - FIND THE API PF THIS" paragraph
# ---------------------------------------------------------------------
$calPara = Find-ParagraphByPrefix $d "Calander"
if ($null -eq $calPara) {
    throw "Could not locate the 'Calander(not fix) ... FIND THE API PF THIS' paragraph"
}
$calRange = $calPara.Range

$calXml = '<w:p ' + $wNs + ' w:rsidR="001A1CB1" w:rsidRDefault="00CB712A">' + `
    '<w:pPr><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:proofErr w:type="gramStart"/>' + `
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Calander</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>(</w:t></w:r>' + `
    '<w:proofErr w:type="gramEnd"/>' + `
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>not fix)</w:t></w:r>' + `
    '<w:r w:rsidR="001A1CB1"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> &#8212; FIND THE API PF THIS</w:t></w:r>' + `
    '</w:p>'

[void]$calRange.InsertXML($calXml)

# ---------------------------------------------------------------------
# Step 2: Bold the "Portfolio" paragraph
# ---------------------------------------------------------------------
$portPara = Find-ParagraphByPrefix $d "Portfolio"
if ($null -eq $portPara) {
    throw "Could not locate the 'Portfolio' paragraph"
}
$portPara.Range.Font.Bold = 1
